$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "67.098.28"
$ws.Range("E2").Value = "  +4.39%  "

$ws.Range("D3").Value2 = "3.255.83"
$ws.Range("E3").Value = "  +2.71%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'578.01"
$ws.Range("E5").Value = "  +2.55%  "

$ws.Range("D6").Value = "'178.00"
$ws.Range("E6").Value = "  +4.48%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "'0.602"
$ws.Range("E8").Value = "  -1.02%  "

$ws.Range("D9").Value2 = "3.254.60"
$ws.Range("E9").Value = "  +2.69%  "

$ws.Range("E10").Value = "  +4.39%  "

$ws.Range("D11").Value = "'6.75"
$ws.Range("E11").Value = "  +1.84%  "

$ws.Range("E12").Value = "  +4.40%  "

$ws.Range("D13").Value2 = "3.821.72"

$ws.Range("E14").Value = "  +0.74%  "

$ws.Range("D15").Value = "'28.12"
$ws.Range("E15").Value = "  +2.73%  "

$ws.Range("D16").Value2 = "67.106.37"
$ws.Range("E16").Value = "  +4.41%  "

$ws.Range("E17").Value = "  +2.99%  "

$ws.Range("D18").Value2 = "3.257.04"
$ws.Range("E18").Value = "  +2.88%  "

$ws.Range("E19").Value = "  +2.19%  "

$ws.Range("D20").Value = "'13.42"
$ws.Range("E20").Value = "  +3.13%  "

$ws.Range("D21").Value = "'373.31"
$ws.Range("E21").Value = "  +5.47%  "

$ws.Range("E22").Value = "  +6.05%  "

$ws.Range("D24").Value = "'71.01"
$ws.Range("E24").Value = "  +2.75%  "

$ws.Range("E25").Value = "  +1.69%  "

$ws.Range("D26").Value2 = "3.397.34"
$ws.Range("E26").Value = "  +2.86%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("D28").Value = "'9.89"
$ws.Range("E28").Value = "  +3.47%  "

$ws.Range("E29").Value = "  +1.96%  "

$ws.Range("E30").Value = "  +0.09%  "

$ws.Range("E31").Value = "  +4.37%  "

$ws.Range("E32").Value = "  +0.22%  "

$ws.Range("D33").Value = "'22.60"
$ws.Range("E33").Value = "  +2.44%  "

$ws.Range("E34").Value = "  -0.07%  "

$ws.Range("E35").Value = "  +5.23%  "

$ws.Range("E36").Value = "  +2.78%  "

$ws.Range("D37").Value = "'166.74"
$ws.Range("E37").Value = "  +7.42%  "

$ws.Range("E38").Value = "  +4.75%  "

$ws.Range("D39").Value = "'0.856"
$ws.Range("E39").Value = "  +5.49%  "

$ws.Range("E40").Value = "  +10.38%  "

$ws.Range("D41").Value = "'27.14"
$ws.Range("E41").Value = "  +4.81%  "

$ws.Range("D42").Value = "'2.59"
$ws.Range("E42").Value = "  +1.74%  "

$ws.Range("D43").Value2 = "2.759.58"
$ws.Range("E43").Value = "  +5.96%  "

$ws.Range("E44").Value = "  +8.30%  "

$ws.Range("D45").Value = "'354.39"
$ws.Range("E45").Value = "  +10.24%  "

$ws.Range("E46").Value = "  +4.95%  "

$ws.Range("D47").Value = "'25.44"
$ws.Range("E47").Value = "  +6.19%  "

$ws.Range("E48").Value = "  +2.10%  "

$ws.Range("E49").Value = "  +2.45%  "

$ws.Range("D50").Value = "'0.0281"
$ws.Range("E50").Value = "  +3.73%  "

$ws.Range("E51").Value = "  +0.63%  "
